# "added product update functionality + a lot of refactoring"
#
# Applies the edits described by the diff:
#  - Row 19 (B19): new wording "Like product (Increase popularity) and
#    identify if the user has liked this product or not", wrap-text style,
#    taller row.
#  - Row 18 (C18): recolor the "Solved" marker cell to the new
#    green-on-green ("done") style instead of the red-on-red style.
#  - New row 32: "Use @Transactional for methods that need more than one
#    Db interaction" entry added to the "Tehnical" table, styled like the
#    other multi-line rows in that table.
#  - Column B widened to fit the longer text.
#  - Selection moved to D11 (matches the author's last cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B widen (45.109375 -> 54.6640625 char-width units) -----------
$ws.Columns.Item(2).ColumnWidth = 53.8333

# --- Row 19: "Like product ..." gets expanded wording --------------------
$newLike = "Like product (Increase popularity) and identify `nif the user has liked this product or not"
$ws.Range("B19").Value = $newLike
$ws.Range("B19").WrapText = $true
$ws.Rows.Item(19).RowHeight = 29.4

# --- Row 18: C18 status cell switches from the red "todo" style to the ---
# --- green "done" style (green font on green fill, same border) ----------
$ws.Range("C30").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Font.Color = 5287936

# --- New row 32: add the @Transactional note to the Tehnical table -------
$ws.Range("A32").Value = 8
$newTx = "Use @Transactional for methods that need more than one `nDb interaction "
$ws.Range("B32").Value = $newTx
$ws.Range("B32").WrapText = $true
$ws.Rows.Item(32).RowHeight = 29.4

# Give C32 the same green fill + medium border used by the other rows in
# this table (copy style from C30 which already carries it).
$ws.Range("C30").Copy() | Out-Null
$ws.Range("C32").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Selection follows the author's last cursor position -----------------
$ws.Range("D11").Select() | Out-Null
